$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1514.7142
$ws.Range("I98").Value = 1267.1666
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 1267.1666
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 230.8334
$ws.Range("N98").Value = -5996

$ws.Range("H122").Value = 1514.7142
$ws.Range("I122").Value = 1267.1666
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3801.4998
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1351.4998
$ws.Range("N122").Value = -13900

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws.Range("H137").Value = 1689.2258
$ws.Range("I137").Value = 1599.125
$ws.Range("J137").Value = 1785.3334
$ws.Range("K137").Value = 4797.375
$ws.Range("L137").Value = 5356.0002
$ws.Range("M137").Value = -2247.375
$ws.Range("N137").Value = -10456.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 491.70834
$ws.Range("I2").Value = 442.78946
$ws.Range("K2").Value = 442.78946
$ws.Range("M2").Value = -329.78946

$ws.Range("H45").Value = 1610.8077
$ws.Range("I45").Value = 1228.1
$ws.Range("J45").Value = 2886.5
$ws.Range("K45").Value = 1228.1
$ws.Range("L45").Value = 2886.5
$ws.Range("M45").Value = -851.0999999999999
$ws.Range("N45").Value = -3640.5

$ws.Range("H116").Value = 491.70834
$ws.Range("I116").Value = 442.78946
$ws.Range("K116").Value = 442.78946
$ws.Range("M116").Value = 1851.21054

$ws.Range("H132").Value = 2751566.5
$ws.Range("I132").Value = 3118.3125
$ws.Range("J132").Value = 6416164
$ws.Range("K132").Value = 9354.9375
$ws.Range("L132").Value = 19248492
$ws.Range("M132").Value = -6824.9375
$ws.Range("N132").Value = -19253552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 491.70834
$ws.Range("I3").Value = 442.78946
$ws.Range("K3").Value = 442.78946
$ws.Range("M3").Value = -328.78946

$ws.Range("H134").Value = 3184.7742
$ws.Range("I134").Value = 3133.3333
$ws.Range("J134").Value = 3292.8
$ws.Range("K134").Value = 9399.999899999999
$ws.Range("L134").Value = 9878.400000000001
$ws.Range("M134").Value = -6864.999899999999
$ws.Range("N134").Value = -14948.4

$ws.Range("H138").Value = 42624
$ws.Range("J138").Value = 50780
$ws.Range("L138").Value = 50780
$ws.Range("N138").Value = -61060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5186.279
$ws.Range("I31").Value = 1404.5264
$ws.Range("J31").Value = 6897.0713
$ws.Range("K31").Value = 1404.5264
$ws.Range("L31").Value = 6897.0713
$ws.Range("M31").Value = -1109.5264
$ws.Range("N31").Value = -7487.0713

$ws.Range("H34").Value = 5186.279
$ws.Range("I34").Value = 1404.5264
$ws.Range("J34").Value = 6897.0713
$ws.Range("K34").Value = 1404.5264
$ws.Range("L34").Value = 6897.0713
$ws.Range("M34").Value = -1202.5264
$ws.Range("N34").Value = -7301.0713

$ws.Range("H58").Value = 2881
$ws.Range("I58").Value = 2735.037
$ws.Range("J58").Value = 3318.889
$ws.Range("K58").Value = 2735.037
$ws.Range("L58").Value = 3318.889
$ws.Range("M58").Value = -2532.037
$ws.Range("N58").Value = -3724.889

$ws.Range("H122").Value = 2006.7
$ws.Range("I122").Value = 2022.3334
$ws.Range("K122").Value = 6067.0002
$ws.Range("M122").Value = -3617.0002

$ws.Range("H132").Value = 2083.0667
$ws.Range("I132").Value = 1769.4445
$ws.Range("J132").Value = 2553.5
$ws.Range("K132").Value = 5308.333500000001
$ws.Range("L132").Value = 7660.5
$ws.Range("M132").Value = -2778.333500000001
$ws.Range("N132").Value = -12720.5

$ws.Range("H134").Value = 10006555
$ws.Range("I134").Value = 15632908
$ws.Range("J134").Value = 4148.3335
$ws.Range("K134").Value = 46898724
$ws.Range("L134").Value = 12445.0005
$ws.Range("M134").Value = -46896189
$ws.Range("N134").Value = -17515.0005

$ws.Range("H135").Value = 45136.5
$ws.Range("J135").Value = 45136.5
$ws.Range("L135").Value = 45136.5
$ws.Range("N135").Value = -55276.5

$ws.Range("H136").Value = 2881
$ws.Range("I136").Value = 2735.037
$ws.Range("J136").Value = 3318.889
$ws.Range("K136").Value = 8205.110999999999
$ws.Range("L136").Value = 9956.667000000001
$ws.Range("M136").Value = -5655.110999999999
$ws.Range("N136").Value = -15056.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 509.27026
$ws.Range("I113").Value = 500.16666
$ws.Range("J113").Value = 517.8946999999999
$ws.Range("K113").Value = 1500.49998
$ws.Range("L113").Value = 1553.6841
$ws.Range("M113").Value = 669.5000199999999
$ws.Range("N113").Value = -5893.6841

$ws.Range("H133").Value = 12790.606
$ws.Range("I133").Value = 9847.777
$ws.Range("J133").Value = 13894.167
$ws.Range("K133").Value = 29543.331
$ws.Range("L133").Value = 41682.501
$ws.Range("M133").Value = -24483.331
$ws.Range("N133").Value = -51802.501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5671.8945
$ws.Range("I7").Value = 5664.4
$ws.Range("J7").Value = 5700
$ws.Range("K7").Value = 5664.4
$ws.Range("L7").Value = 5700
$ws.Range("M7").Value = -5552.4
$ws.Range("N7").Value = -5924

$ws.Range("H16").Value = 10989775
$ws.Range("I16").Value = 721.36365
$ws.Range("J16").Value = 71429570
$ws.Range("K16").Value = 721.36365
$ws.Range("L16").Value = 71429570
$ws.Range("M16").Value = -551.36365
$ws.Range("N16").Value = -71429910

$ws.Range("H93").Value = 12634
$ws.Range("I93").Value = 13588.25
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 13588.25
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -12340.25
$ws.Range("N93").Value = -7496

$ws.Range("H122").Value = 5450.909
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 6653.3335
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 19960.0005
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -24860.0005

$ws.Range("H126").Value = 5671.8945
$ws.Range("I126").Value = 5664.4
$ws.Range("J126").Value = 5700
$ws.Range("K126").Value = 16993.2
$ws.Range("L126").Value = 17100
$ws.Range("M126").Value = -14523.2
$ws.Range("N126").Value = -22040

$ws.Range("H137").Value = 35312.5
$ws.Range("I137").Value = 23000
$ws.Range("J137").Value = 37071.43
$ws.Range("K137").Value = 23000
$ws.Range("L137").Value = 37071.43
$ws.Range("N137").Value = -47271.43
$ws.Range("M137").Value = -17900

$ws.Range("H139").Value = 4862350
$ws.Range("J139").Value = 49800.332
$ws.Range("L139").Value = 49800.332
$ws.Range("N139").Value = -60080.332

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 34655
$ws.Range("J137").Value = 34655
$ws.Range("L137").Value = 34655
$ws.Range("N137").Value = -44855
